$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 250
$ws1.Range("F3").Value = 444
$ws1.Range("F6").Value = 543
$ws1.Range("F8").Value = 28
$ws1.Range("F9").Value = 267
$ws1.Range("F10").Value = 380
$ws1.Range("F11").Value = 348
$ws1.Range("F12").Value = 643
$ws1.Range("F13").Value = 749
$ws1.Range("F14").Value = 1502
$ws1.Range("F15").Value = 1502
$ws1.Range("F16").Value = 882
$ws1.Range("F17").Value = 27
$ws1.Range("F18").Value = 1349
$ws1.Range("F20").Value = 287
$ws1.Range("F23").Value = 97
$ws1.Range("F24").Value = 6549
$ws1.Range("F25").Value = 4883
$ws1.Range("F26").Value = 141
$ws1.Range("F29").Value = 151
$ws1.Range("F32").Value = 1274
$ws1.Range("F33").Value = 191
$ws1.Range("F34").Value = 243
$ws1.Range("F35").Value = 599
$ws1.Range("F36").Value = 18
$ws1.Range("F38").Value = 237
$ws1.Range("F39").Value = 149
$ws1.Range("F41").Value = 61
$ws1.Range("F43").Value = 94

# Sheet 3: "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 190
$ws3.Range("F5").Value = 48

# Sheet 4: "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 250
$ws4.Range("F4").Value = 444
$ws4.Range("F7").Value = 190
$ws4.Range("F8").Value = 48
$ws4.Range("F10").Value = 543
$ws4.Range("F12").Value = 28
$ws4.Range("F13").Value = 267
$ws4.Range("F15").Value = 380
$ws4.Range("F16").Value = 348
$ws4.Range("F17").Value = 643
$ws4.Range("F18").Value = 749
$ws4.Range("F19").Value = 1502
$ws4.Range("F20").Value = 1502
$ws4.Range("F21").Value = 882
$ws4.Range("F22").Value = 27
$ws4.Range("F23").Value = 1349
$ws4.Range("F25").Value = 287
$ws4.Range("F27").Value = 97
$ws4.Range("F30").Value = 6549
$ws4.Range("F31").Value = 4883
$ws4.Range("F32").Value = 141
$ws4.Range("F34").Value = 1274
$ws4.Range("F35").Value = 191
$ws4.Range("F36").Value = 243
$ws4.Range("F38").Value = 599
$ws4.Range("F40").Value = 18
$ws4.Range("F44").Value = 237
$ws4.Range("F46").Value = 61
$ws4.Range("F48").Value = 94

Write-Output "Done applying updates"
